$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that are no longer needed (original row numbers 8, 6, 5, 3)
# Delete from bottom to top so row indices of earlier rows stay valid.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(3).Delete()

# Update the selection to match the saved view state
$ws.Range("B7").Select()
